$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.409.36"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.841.89"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "259.07"
$ws.Range("E5").Value = "  -7.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5098"
$ws.Range("E7").Value = "  -0.61%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3213"
$ws.Range("E8").Value = "  -8.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06748"
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.84"
$ws.Range("E10").Value = "  -5.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7682"
$ws.Range("E11").Value = "  -5.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07678"
$ws.Range("D13").Value = "1.875.10"
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.44"
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.014"
$ws.Range("E15").Value = "  -1.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("E17").Value = "  -1.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007887"
$ws.Range("E19").Value = "  -2.66%  "
$ws.Range("D20").Value = "26.453.57"
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").Value = "2.084.72"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.575"
$ws.Range("E22").Value = "  -4.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.520"
$ws.Range("E23").Value = "  -5.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.951"
$ws.Range("E24").Value = "  -4.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.326"
$ws.Range("E25").Value = "  -1.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "145.12"
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.662"
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.95"
$ws.Range("E28").Value = "  -2.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "110.87"
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("E30").Value = "  -4.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.148"
$ws.Range("E31").Value = "  -4.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08701"
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04812"
$ws.Range("E33").Value = "  -1.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.131"
$ws.Range("E34").Value = "  -3.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.845"
$ws.Range("E35").Value = "  -0.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6839"
$ws.Range("E36").Value = "  -7.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.082"
$ws.Range("E37").Value = "  -5.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01801"
$ws.Range("E38").Value = "  -3.31%  "
$ws.Range("E39").Value = "  -7.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4893"
$ws.Range("E40").Value = "  -6.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "113.09"
$ws.Range("E41").Value = "  -2.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9028"
$ws.Range("E42").Value = "  -6.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.102"
$ws.Range("E43").Value = "  -2.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.758"
$ws.Range("E45").Value = "  -3.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4236"
$ws.Range("E46").Value = "  -6.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1270"
$ws.Range("E47").Value = "  -6.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.126"
$ws.Range("E48").Value = "  -2.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05892"
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.94"
$ws.Range("E50").Value = "  -3.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.419"
$ws.Range("E51").Value = "  -5.71%  "
